$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("company_list")

# Row 2
$ws.Range("D2").Value = 15201
$ws.Range("E2").Value = 230
$ws.Range("F2").Value = -797
$ws.Range("G2").Value = -1285
$ws.Range("H2").Value = -1350
$ws.Range("I2").Value = -1214
$ws.Range("J2").Value = -137
$ws.Range("K2").Value = 27222
$ws.Range("L2").Value = 16149
$ws.Range("M2").Value = 11073
$ws.Range("N2").Value = 8985
$ws.Range("O2").Value = 2088
$ws.Range("P2").Value = 1551
$ws.Range("Q2").Value = 502
$ws.Range("R2").Value = -1693
$ws.Range("S2").Value = 1330
$ws.Range("T2").Value = 1172
$ws.Range("U2").Value = -670
$ws.Range("V2").Value = 9908
$ws.Range("W2").Value = 1.51
$ws.Range("X2").Value = -8.880000000000001
$ws.Range("Y2").Value = -12.71
$ws.Range("Z2").Value = -5.14
$ws.Range("AA2").Value = 145.84
$ws.Range("AB2").Value = 476.93
$ws.Range("AC2").Value = -4110
$ws.Range("AD2").Value = -1.98
$ws.Range("AE2").Value = 34439
$ws.Range("AF2").Value = 0.24
$ws.Range("AG2").Value = 200
$ws.Range("AH2").Value = 2.46
$ws.Range("AI2").Value = -2.04
$ws.Range("AJ2").Value = 29529812

# Row 3
$ws.Range("D3").Value = 2752
$ws.Range("E3").Value = 90
$ws.Range("F3").Value = -683
$ws.Range("G3").Value = -1200
$ws.Range("H3").Value = -1643
$ws.Range("I3").Value = -1379
$ws.Range("J3").Value = -264
$ws.Range("K3").Value = 25858
$ws.Range("L3").Value = 16487
$ws.Range("M3").Value = 9371
$ws.Range("N3").Value = 7565
$ws.Range("O3").Value = 1806
$ws.Range("P3").Value = 1551
$ws.Range("Q3").Value = 588
$ws.Range("R3").Value = -679
$ws.Range("S3").Value = 568
$ws.Range("T3").Value = 437
$ws.Range("U3").Value = 150
$ws.Range("V3").Value = 10581
$ws.Range("W3").Value = 3.27
$ws.Range("X3").Value = -59.71
$ws.Range("Y3").Value = -16.67
$ws.Range("Z3").Value = -6.19
$ws.Range("AA3").Value = 175.94
$ws.Range("AB3").Value = 373.26
$ws.Range("AC3").Value = -4671
$ws.Range("AD3").Value = -1.29
$ws.Range("AE3").Value = 28996
$ws.Range("AF3").Value = 0.21
$ws.Range("AG3").Value = 0
$ws.Range("AH3").Value = 0
$ws.Range("AI3").Value = 0
$ws.Range("AJ3").Value = 29529812

# Row 4
$ws.Range("D4").Value = 9683
$ws.Range("E4").Value = -841
$ws.Range("F4").Value = -978
$ws.Range("G4").Value = -3052
$ws.Range("H4").Value = -2911
$ws.Range("I4").Value = -2668
$ws.Range("J4").Value = -243
$ws.Range("K4").Value = 22775
$ws.Range("L4").Value = 16225
$ws.Range("M4").Value = 6550
$ws.Range("N4").Value = 4994
$ws.Range("O4").Value = 1556
$ws.Range("P4").Value = 1551
$ws.Range("Q4").Value = 269
$ws.Range("R4").Value = 65
$ws.Range("S4").Value = -433
$ws.Range("T4").Value = 158
$ws.Range("U4").Value = 112
$ws.Range("V4").Value = 10287
$ws.Range("W4").Value = -8.69
$ws.Range("X4").Value = -30.06
$ws.Range("Y4").Value = -42.5
$ws.Range("Z4").Value = -11.97
$ws.Range("AA4").Value = 247.72
$ws.Range("AB4").Value = 202.27
$ws.Range("AC4").Value = -9036
$ws.Range("AD4").Value = -0.5600000000000001
$ws.Range("AE4").Value = 19140
$ws.Range("AF4").Value = 0.26
$ws.Range("AG4").Value = 0
$ws.Range("AH4").Value = 0
$ws.Range("AI4").Value = 0
$ws.Range("AJ4").Value = 29529812

# Row 5
$ws.Range("D5").Value = 9067
$ws.Range("E5").Value = 243
$ws.Range("F5").Value = -618
$ws.Range("G5").Value = -1126
$ws.Range("H5").Value = -1596
$ws.Range("I5").Value = -1394
$ws.Range("J5").Value = -202
$ws.Range("K5").Value = 19186
$ws.Range("L5").Value = 14897
$ws.Range("M5").Value = 4290
$ws.Range("N5").Value = 3313
$ws.Range("O5").Value = 977
$ws.Range("P5").Value = 1551
$ws.Range("Q5").Value = 350
$ws.Range("R5").Value = 341
$ws.Range("S5").Value = -772
$ws.Range("T5").Value = 119
$ws.Range("U5").Value = 232
$ws.Range("V5").Value = 9244
$ws.Range("W5").Value = 2.68
$ws.Range("X5").Value = -17.61
$ws.Range("Y5").Value = -33.57
$ws.Range("Z5").Value = -7.61
$ws.Range("AA5").Value = 347.26
$ws.Range("AB5").Value = 113.94
$ws.Range("AC5").Value = -4722
$ws.Range("AD5").Value = -0.85
$ws.Range("AE5").Value = 12698
$ws.Range("AF5").Value = 0.32
$ws.Range("AG5").Value = 0
$ws.Range("AH5").Value = 0
$ws.Range("AI5").Value = 0
$ws.Range("AJ5").Value = 29529812

# Row 6
$ws.Range("D6").Value = 9460
$ws.Range("E6").Value = 324
$ws.Range("F6").Value = 277
$ws.Range("G6").Value = 404
$ws.Range("H6").Value = -3
$ws.Range("I6").Value = 90
$ws.Range("K6").Value = 18736
$ws.Range("L6").Value = 15066
$ws.Range("M6").Value = 3670
$ws.Range("N6").Value = 3403
$ws.Range("P6").Value = 1551
$ws.Range("Q6").Value = 572
$ws.Range("R6").Value = -590
$ws.Range("S6").Value = -566
$ws.Range("T6").Value = 350
$ws.Range("U6").Value = 221
$ws.Range("V6").Value = 9311
$ws.Range("W6").Value = 3.42
$ws.Range("X6").Value = -0.03
$ws.Range("Y6").Value = 2.69
$ws.Range("Z6").Value = -0.01
$ws.Range("AA6").Value = 410.45
$ws.Range("AB6").Value = 116.58
$ws.Range("AC6").Value = 306
$ws.Range("AD6").Value = 9.57
$ws.Range("AE6").Value = 13043
$ws.Range("AF6").Value = 0.22
$ws.Range("AJ6").Value = 29529812
$ws.Range("AG6").ClearContents()
$ws.Range("AH6").ClearContents()
$ws.Range("AI6").ClearContents()

# Row 7
$ws.Range("D7").ClearContents()
$ws.Range("E7").ClearContents()
$ws.Range("G7").ClearContents()
$ws.Range("H7").ClearContents()
$ws.Range("I7").ClearContents()
$ws.Range("K7").ClearContents()
$ws.Range("L7").ClearContents()
$ws.Range("M7").ClearContents()
$ws.Range("N7").ClearContents()
$ws.Range("P7").ClearContents()
$ws.Range("Q7").ClearContents()
$ws.Range("R7").ClearContents()
$ws.Range("S7").ClearContents()
$ws.Range("T7").ClearContents()
$ws.Range("U7").ClearContents()
$ws.Range("W7").ClearContents()
$ws.Range("X7").ClearContents()
$ws.Range("Y7").ClearContents()
$ws.Range("Z7").ClearContents()
$ws.Range("AA7").ClearContents()
$ws.Range("AC7").ClearContents()
$ws.Range("AD7").ClearContents()
$ws.Range("AE7").ClearContents()
$ws.Range("AF7").ClearContents()
$ws.Range("AG7").ClearContents()
$ws.Range("AH7").ClearContents()
$ws.Range("AI7").ClearContents()

# Row 8
$ws.Range("D8").ClearContents()
$ws.Range("E8").ClearContents()
$ws.Range("G8").ClearContents()
$ws.Range("H8").ClearContents()
$ws.Range("I8").ClearContents()
$ws.Range("K8").ClearContents()
$ws.Range("L8").ClearContents()
$ws.Range("M8").ClearContents()
$ws.Range("N8").ClearContents()
$ws.Range("P8").ClearContents()
$ws.Range("Q8").ClearContents()
$ws.Range("R8").ClearContents()
$ws.Range("S8").ClearContents()
$ws.Range("T8").ClearContents()
$ws.Range("U8").ClearContents()
$ws.Range("W8").ClearContents()
$ws.Range("X8").ClearContents()
$ws.Range("Y8").ClearContents()
$ws.Range("Z8").ClearContents()
$ws.Range("AA8").ClearContents()
$ws.Range("AC8").ClearContents()
$ws.Range("AD8").ClearContents()
$ws.Range("AE8").ClearContents()
$ws.Range("AF8").ClearContents()
$ws.Range("AG8").ClearContents()
$ws.Range("AH8").ClearContents()
$ws.Range("AI8").ClearContents()

# Row 9
$ws.Range("D9").ClearContents()
$ws.Range("E9").ClearContents()
$ws.Range("G9").ClearContents()
$ws.Range("H9").ClearContents()
$ws.Range("I9").ClearContents()
$ws.Range("K9").ClearContents()
$ws.Range("L9").ClearContents()
$ws.Range("M9").ClearContents()
$ws.Range("N9").ClearContents()
$ws.Range("P9").ClearContents()
$ws.Range("Q9").ClearContents()
$ws.Range("R9").ClearContents()
$ws.Range("S9").ClearContents()
$ws.Range("T9").ClearContents()
$ws.Range("U9").ClearContents()
$ws.Range("W9").ClearContents()
$ws.Range("X9").ClearContents()
$ws.Range("Y9").ClearContents()
$ws.Range("Z9").ClearContents()
$ws.Range("AA9").ClearContents()
$ws.Range("AC9").ClearContents()
$ws.Range("AD9").ClearContents()
$ws.Range("AE9").ClearContents()
$ws.Range("AF9").ClearContents()
$ws.Range("AG9").ClearContents()
$ws.Range("AH9").ClearContents()
$ws.Range("AI9").ClearContents()
